# Update "想去人数" (F column) values across the exhibition/show/all-type sheets,
# reflecting refreshed data as of the regeneration run (gh-pages output).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 20525
    3  = 805
    4  = 323
    6  = 19
    7  = 7690
    8  = 529
    10 = 287
    12 = 172
    13 = 137
    15 = 240
    16 = 205
    17 = 1350
    18 = 468
    20 = 693
    24 = 333
    25 = 1142
    26 = 40
    27 = 26
    28 = 194
    29 = 5213
    30 = 580
    31 = 97
    32 = 4914
    33 = 28
    34 = 94
    36 = 12802
    38 = 96
    39 = 36
    40 = 62
    41 = 284
    42 = 395
    43 = 4022
    44 = 325
    45 = 99
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    2 = 228
    4 = 35
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Cells.Item($row, 6).Value = $sheet2Updates[$row]
}

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 20526
    3  = 805
    4  = 323
    6  = 20
    7  = 7690
    8  = 529
    10 = 287
    12 = 172
    13 = 137
    15 = 240
    16 = 205
    17 = 1350
    18 = 468
    20 = 693
    24 = 333
    25 = 1142
    26 = 40
    27 = 26
    28 = 194
    29 = 228
    30 = 5213
    31 = 580
    33 = 97
    34 = 35
    35 = 4914
    36 = 28
    37 = 94
    39 = 12802
    41 = 96
    42 = 36
    43 = 62
    44 = 284
    45 = 395
    46 = 4022
    47 = 325
    48 = 99
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
